# Auto-generated edit script: updates cryptos list values (price/volume)
# matching the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.417.60'
$ws.Range('E2').Value = '  -1.94%  '
$ws.Range('D3').Value = '3.487.98'
$ws.Range('E3').Value = '  -2.32%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '''609.64'
$ws.Range('E5').Value = '  +4.66%  '
$ws.Range('D6').Value = '''186.09'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('E9').Value = '  -3.95%  '
$ws.Range('D10').Value = '''0.647'
$ws.Range('E10').Value = '  -0.71%  '
$ws.Range('D11').Value = '''52.97'
$ws.Range('E11').Value = '  -2.87%  '
$ws.Range('E12').Value = '  -3.77%  '
$ws.Range('D13').Value = '''9.53'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').Value = '4.040.25'
$ws.Range('E14').Value = '  -2.35%  '
$ws.Range('D15').Value = '''602.44'
$ws.Range('E15').Value = '  +4.82%  '
$ws.Range('D16').Value = '69.470.12'
$ws.Range('E16').Value = '  -1.89%  '
$ws.Range('D17').Value = '''18.82'
$ws.Range('E17').Value = '  -2.17%  '
$ws.Range('D18').Value = '''12.55'
$ws.Range('E18').Value = '  -1.56%  '
$ws.Range('D19').Value = '3.466.18'
$ws.Range('E19').Value = '  -3.31%  '
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').Value = '''0.985'
$ws.Range('E21').Value = '  -1.80%  '
$ws.Range('D22').Value = '''17.12'
$ws.Range('E22').Value = '  -2.68%  '
$ws.Range('E23').Value = '  +12.54%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '''5.09'
$ws.Range('E24').Value = '  +3.96%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '''4.63'
$ws.Range('E25').Value = '  +1.52%  '
$ws.Range('D26').Value = '''3.05'
$ws.Range('E26').Value = '  +2.98%  '
$ws.Range('E27').Value = '  -2.68%  '
$ws.Range('E28').Value = '  +4.92%  '
$ws.Range('D29').Value = '''33.52'
$ws.Range('E29').Value = '  +3.26%  '
$ws.Range('E30').Value = '  -3.39%  '
$ws.Range('D31').Value = '''4.10'
$ws.Range('E31').Value = '  +16.24%  '
$ws.Range('D32').Value = '''12.41'
$ws.Range('E32').Value = '  +0.90%  '
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').Value = '''63.30'
$ws.Range('E34').Value = '  +0.35%  '
$ws.Range('D35').Value = '''3.18'
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('D37').Value = '''521.44'
$ws.Range('E37').Value = '  -5.53%  '
$ws.Range('E38').Value = '  +6.91%  '
$ws.Range('D39').Value = '3.611.41'
$ws.Range('E39').Value = '  +0.65%  '
$ws.Range('E41').Value = '  -3.33%  '
$ws.Range('D42').Value = '0.0₃0778'
$ws.Range('E42').Value = '  -3.72%  '
$ws.Range('D43').Value = '''0.136'
$ws.Range('E43').Value = '  -0.94%  '
$ws.Range('D44').Value = '''0.0461'
$ws.Range('E44').Value = '  -0.65%  '
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('E46').Value = '  +3.10%  '
$ws.Range('D47').Value = '''3.34'
$ws.Range('E47').Value = '  -3.92%  '
$ws.Range('E48').Value = '  -5.79%  '
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('D50').Value = '''1.36'
$ws.Range('E50').Value = '  -10.45%  '
$ws.Range('E51').Value = '  -7.87%  '
